# Applies the "PO Forecast" update:
#  1. Rename the "Requested quantity" header on "Weekly Quantity" to "Weekly_PO_Qty"
#  2. Rename the "Requested quantity" header on "Monthly Trend" to "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename existing headers -----------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the "PO Forecast" sheet after the last existing sheet -----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# header row
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$rows = @(
    @(45242.99999999999, 124, 124.0008541478623, 124.0008541546977),
    @(45277.99999999999, 4, 4.000846082278708, 4.000846088981214),
    @(45284.99999999999, 0, -19.99915553236481, -19.99915551617428),
    @(45291.99999999999, 0, -43.99915726422083, -43.99915698565845),
    @(45298.99999999999, 0, -67.99915903650546, -67.99915842212327),
    @(45305.99999999999, 0, -91.99916091174501, -91.99915970071878),
    @(45312.99999999999, 0, -115.999162844506, -115.9991609465024),
    @(45319.99999999999, 0, -139.9991647514637, -139.9991622297664),
    @(45326.99999999999, 0, -163.9991666942067, -163.999163399426),
    @(45333.99999999999, 0, -187.9991687052836, -187.9991645904078)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r,1).Value = $row[0]
    $wsForecast.Cells.Item($r,2).Value = $row[1]
    $wsForecast.Cells.Item($r,3).Value = $row[2]
    $wsForecast.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

# --- formatting: mirror the header / date styles used on the other sheets ---
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

$excel.CutCopyMode = $false
